# Insert a new weekly record row at row 611 (shifting the existing
# rows 611-635 down to 612-636) and populate the new row with the
# latest week's data for Femacal de La Calera - Ajo.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 611..635 down by one, inserting a fresh row at 611.
$ws.Rows.Item(611).EntireRow.Insert()

# Populate the newly inserted row 611 with the new weekly record.
$ws.Cells.Item(611, 1).Value = 3
$ws.Cells.Item(611, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(611, 3).Value = "Coquimbo"
$ws.Cells.Item(611, 4).Value = 44939
$ws.Cells.Item(611, 5).Value = 5
$ws.Cells.Item(611, 6).Value = 100112003
$ws.Cells.Item(611, 7).Value = "Ajo"
$ws.Cells.Item(611, 8).Value = "Chino"
$ws.Cells.Item(611, 9).Value = "Primera"
$ws.Cells.Item(611, 10).Value = 85
$ws.Cells.Item(611, 11).Value = 15500
$ws.Cells.Item(611, 12).Value = 16000
$ws.Cells.Item(611, 13).Value = 15765
$ws.Cells.Item(611, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(611, 15).Value = "China"
$ws.Cells.Item(611, 16).Value = 1576
$ws.Cells.Item(611, 17).Value = 10
$ws.Cells.Item(611, 18).Value = "Hortaliza"
